$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "83÷3="; New = "58÷5=" },
    @{ Old = "92÷5="; New = "73÷8=" },
    @{ Old = "30÷4="; New = "31÷7=" },
    @{ Old = "10÷7="; New = "52÷3=" },
    @{ Old = "42÷6="; New = "48÷9=" },
    @{ Old = "66÷3="; New = "78÷8=" },
    @{ Old = "21÷4="; New = "73÷8=" },
    @{ Old = "73÷7="; New = "43÷3=" },
    @{ Old = "93÷9="; New = "93÷8=" },
    @{ Old = "97÷5="; New = "73÷4=" },
    @{ Old = "59÷3="; New = "58÷2=" },
    @{ Old = "11÷2="; New = "86÷3=" },
    @{ Old = "87÷2="; New = "39÷2=" },
    @{ Old = "95÷7="; New = "70÷2=" },
    @{ Old = "98÷9="; New = "58÷9=" },
    @{ Old = "24÷5="; New = "61÷3=" },
    @{ Old = "91÷8="; New = "95÷4=" },
    @{ Old = "17÷8="; New = "33÷2=" },
    @{ Old = "55÷5="; New = "37÷3=" },
    @{ Old = "41÷2="; New = "78÷9=" },
    @{ Old = "11÷8="; New = "88÷7=" },
    @{ Old = "77÷2="; New = "39÷2=" },
    @{ Old = "63÷4="; New = "58÷2=" },
    @{ Old = "12÷7="; New = "14÷3=" },
    @{ Old = "44÷4="; New = "56÷3=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
